# B6-PowerPoint.pptx edit: 2020-05-23
#
# 1) Re-style the three tables (slides 14, 15, 16) from the custom
#    "Table_0" style {77EF1E9C-194D-4E96-BE27-9A32FCF7CE7B} to the
#    built-in table style {94CB7813-182D-474B-8193-B109E963B643}.
# 2) Re-colour the deck's theme (theme1.xml, the one the slide master
#    actually uses) from the "Integral"/"Red Violet" palette to the
#    stock Office palette, matching the target clrScheme:
#      dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#      accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#      accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newTableStyleId = "{94CB7813-182D-474B-8193-B109E963B643}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Theme colours --------------------------------------------------
# RGB() below is PowerPoint's little-endian 0xBBGGRR long, built from the
# target hex triples so the written <a:srgbClr val="..."/> match exactly.
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$themeColorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColorScheme.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [math]::Floor($hex / 0x10000) % 0x100
    $g = [math]::Floor($hex / 0x100) % 0x100
    $b = $hex % 0x100
    $themeColorScheme.Colors($i).RGB = $r + ($g * 0x100) + ($b * 0x10000)
}
